$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 681
$ws.Range("I6").Value = 351.25
$ws.Range("K6").Value = 1053.75
$ws.Range("M6").Value = -941.75
$ws.Range("H19").Value = 1749.75
$ws.Range("I19").Value = 1199
$ws.Range("J19").Value = 1933.3334
$ws.Range("K19").Value = 1199
$ws.Range("L19").Value = 1933.3334
$ws.Range("M19").Value = -1024
$ws.Range("N19").Value = -2283.3334
$ws.Range("H33").Value = 432.66666
$ws.Range("I33").Value = 361.33334
$ws.Range("J33").Value = 646.6667
$ws.Range("K33").Value = 361.33334
$ws.Range("L33").Value = 646.6667
$ws.Range("M33").Value = -132.33334
$ws.Range("N33").Value = -1104.6667
$ws.Range("H46").Value = 3259
$ws.Range("I46").Value = 3648
$ws.Range("J46").Value = 2999.6667
$ws.Range("K46").Value = 10944
$ws.Range("L46").Value = 8999.000100000001
$ws.Range("M46").Value = -10825
$ws.Range("N46").Value = -9237.000100000001
$ws.Range("H60").Value = 3259
$ws.Range("I60").Value = 3648
$ws.Range("J60").Value = 2999.6667
$ws.Range("K60").Value = 10944
$ws.Range("L60").Value = 8999.000100000001
$ws.Range("M60").Value = -10460
$ws.Range("N60").Value = -9967.000100000001
$ws.Range("H64").Value = 4400
$ws.Range("J64").Value = 6000
$ws.Range("L64").Value = 6000
$ws.Range("N64").Value = -6496
$ws.Range("H67").Value = 4400
$ws.Range("J67").Value = 6000
$ws.Range("L67").Value = 6000
$ws.Range("N67").Value = -7716
$ws.Range("H69").Value = 7580.2354
$ws.Range("I69").Value = 4049.6667
$ws.Range("J69").Value = 8336.786
$ws.Range("K69").Value = 12149.0001
$ws.Range("L69").Value = 25010.358
$ws.Range("M69").Value = -11275.0001
$ws.Range("N69").Value = -26758.358
$ws.Range("H72").Value = 7580.2354
$ws.Range("I72").Value = 4049.6667
$ws.Range("J72").Value = 8336.786
$ws.Range("K72").Value = 36447.0003
$ws.Range("L72").Value = 75031.07399999999
$ws.Range("M72").Value = -32079.0003
$ws.Range("N72").Value = -83767.07399999999
$ws.Range("H74").Value = 14749.75
$ws.Range("I74").Value = 13331.667
$ws.Range("K74").Value = 13331.667
$ws.Range("M74").Value = -12395.667
$ws.Range("H77").Value = 14749.75
$ws.Range("I77").Value = 13331.667
$ws.Range("K77").Value = 66658.33499999999
$ws.Range("M77").Value = -61978.33499999999
$ws.Range("H103").Value = 1136.2858
$ws.Range("J103").Value = 2542
$ws.Range("L103").Value = 7626
$ws.Range("N103").Value = -8798

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5912.2104
$ws.Range("I61").Value = 6156.8
$ws.Range("K61").Value = 6156.8
$ws.Range("M61").Value = -5944.8
$ws.Range("H136").Value = 5912.2104
$ws.Range("I136").Value = 6156.8
$ws.Range("K136").Value = 18470.4
$ws.Range("M136").Value = -15920.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6006.826
$ws.Range("I86").Value = 5407.5
$ws.Range("K86").Value = 5407.5
$ws.Range("M86").Value = -4284.5
$ws.Range("H89").Value = 6006.826
$ws.Range("I89").Value = 5407.5
$ws.Range("K89").Value = 27037.5
$ws.Range("M89").Value = -21421.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5022
$ws.Range("I22").Value = 1948.1428
$ws.Range("J22").Value = 8608.166999999999
$ws.Range("K22").Value = 1948.1428
$ws.Range("L22").Value = 8608.166999999999
$ws.Range("M22").Value = -1598.1428
$ws.Range("N22").Value = -9308.166999999999
$ws.Range("H31").Value = 33662.793
$ws.Range("I31").Value = 2610.7083
$ws.Range("K31").Value = 2610.7083
$ws.Range("M31").Value = -2315.7083
$ws.Range("H34").Value = 33662.793
$ws.Range("I34").Value = 2610.7083
$ws.Range("K34").Value = 2610.7083
$ws.Range("M34").Value = -2408.7083
$ws.Range("H134").Value = 1982.7812
$ws.Range("I134").Value = 1183.091
$ws.Range("K134").Value = 3549.273
$ws.Range("M134").Value = -1014.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 167163.33
$ws.Range("I7").Value = 745
$ws.Range("J7").Value = 500000
$ws.Range("K7").Value = 2235
$ws.Range("L7").Value = 1500000
$ws.Range("M7").Value = -2123
$ws.Range("N7").Value = -1500224
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H75").Value = 250008000
$ws.Range("I75").Value = 1000000000
$ws.Range("J75").Value = 10671.333
$ws.Range("K75").Value = 3000000000
$ws.Range("L75").Value = 32013.999
$ws.Range("M75").Value = -2999999002
$ws.Range("N75").Value = -34009.999
$ws.Range("H78").Value = 250008000
$ws.Range("I78").Value = 1000000000
$ws.Range("J78").Value = 10671.333
$ws.Range("K78").Value = 9000000000
$ws.Range("L78").Value = 96041.997
$ws.Range("M78").Value = -8999995008
$ws.Range("N78").Value = -106025.997
$ws.Range("H80").Value = 10399
$ws.Range("J80").Value = 10399
$ws.Range("L80").Value = 31197
$ws.Range("N80").Value = -33069
$ws.Range("H83").Value = 10399
$ws.Range("J83").Value = 10399
$ws.Range("L83").Value = 93591
$ws.Range("N83").Value = -102951
$ws.Range("H126").Value = 4088.3333
$ws.Range("J126").Value = 4925
$ws.Range("L126").Value = 14775
$ws.Range("N126").Value = -24655

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 380.06668
$ws.Range("I107").Value = 397.2857
$ws.Range("K107").Value = 397.2857
$ws.Range("M107").Value = 1522.7143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3566.087
$ws.Range("I61").Value = 2765.3157
$ws.Range("J61").Value = 7369.75
$ws.Range("K61").Value = 2765.3157
$ws.Range("L61").Value = 7369.75
$ws.Range("M61").Value = -2563.3157
$ws.Range("N61").Value = -7773.75
$ws.Range("H69").Value = 36999
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 36999
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H113").Value = 3566.087
$ws.Range("I113").Value = 2765.3157
$ws.Range("J113").Value = 7369.75
$ws.Range("K113").Value = 2765.3157
$ws.Range("L113").Value = 7369.75
$ws.Range("M113").Value = -595.3157000000001
$ws.Range("N113").Value = -11709.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2383.3157
$ws.Range("J136").Value = 5974
$ws.Range("L136").Value = 17922
$ws.Range("N136").Value = -23022
